# refactor currency conversion, now explicit source and target amounts
#
# currency_conversions header goes from:
#   date | foreign_amount | source_fees | source_currency | target_currency | comment
# to:
#   date | source_amount | source_fees | source_currency | target_amount | target_fees | target_currency | comment
#
# i.e. "foreign_amount" is renamed to "source_amount", and two brand new
# columns - "target_amount" and "target_fees" - are inserted right before the
# existing "target_currency" column (pushing it, and "comment" after it, two
# slots to the right).

$xlShiftToRight = -4161

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_conversions")

# Insert two new header cells just before target_currency (column E), copying
# the format of an existing bold header cell (source_currency, column D) so
# the new cells come out with the same style as their neighbours.
$ws.Range("D1").Copy()
$ws.Range("E1").Insert($xlShiftToRight)
$ws.Range("D1").Copy()
$ws.Range("E1").Insert($xlShiftToRight)

# foreign_amount -> source_amount
$ws.Cells.Item(1, 2).Value = "source_amount"

# the two newly inserted cells become target_amount / target_fees;
# target_currency and comment have already shifted right to G1 / H1.
$ws.Cells.Item(1, 5).Value = "target_amount"
$ws.Cells.Item(1, 6).Value = "target_fees"

# currency_conversions becomes the selected/active sheet (was money_transfers).
$ws.Activate()
